$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date values in columns F and G for rows 2-5 (shift by +1 year / +365 days)
$ws.Range("F2").Value = 42297
$ws.Range("G2").Value = 42556

$ws.Range("F3").Value = 42297
$ws.Range("G3").Value = 42556

$ws.Range("F4").Value = 42297
$ws.Range("G4").Value = 42556

$ws.Range("F5").Value = 42297
$ws.Range("G5").Value = 42556

# Update the active selection/cell shown on the sheet view
$ws.Range("G6").Select()
